$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 426
$ws.Range("I2").Value = 281.8
$ws.Range("K2").Value = 281.8
$ws.Range("M2").Value = -168.8

# Row 29
$ws.Range("H29").Value = 91.25
$ws.Range("I29").Value = 91.25
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 273.75
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 7.25
$ws.Range("N29").ClearContents()

# Row 40
$ws.Range("H40").Value = 2376.7646
$ws.Range("I40").Value = 2977.2222
$ws.Range("J40").Value = 1701.25
$ws.Range("K40").Value = 2977.2222
$ws.Range("L40").Value = 1701.25
$ws.Range("M40").Value = -2802.2222
$ws.Range("N40").Value = -2051.25

# Row 93
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -39992

# Row 98
$ws.Range("H98").Value = 2058.1516
$ws.Range("I98").Value = 2084.9688
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 2084.9688
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = -586.9688000000001
$ws.Range("N98").Value = -4196

# Row 122
$ws.Range("H122").Value = 2058.1516
$ws.Range("I122").Value = 2084.9688
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 6254.9064
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -3804.9064
$ws.Range("N122").Value = -8500

# Row 137
$ws.Range("H137").Value = 2684.2354
$ws.Range("I137").Value = 1503.2391
$ws.Range("J137").Value = 5153.591
$ws.Range("K137").Value = 4509.7173
$ws.Range("L137").Value = 15460.773
$ws.Range("M137").Value = -1959.7173
$ws.Range("N137").Value = -20560.773

# Row 141
$ws.Range("H141").Value = 2741.8333
$ws.Range("I141").Value = 2118.25
$ws.Range("J141").Value = 3989
$ws.Range("K141").Value = 6354.75
$ws.Range("L141").Value = 11967
$ws.Range("M141").Value = -1174.75
$ws.Range("N141").Value = -22327

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 19473.75
$ws.Range("I3").Value = 19800
$ws.Range("J3").Value = 19427.143
$ws.Range("K3").Value = 19800
$ws.Range("L3").Value = 19427.143
$ws.Range("M3").Value = -19685
$ws.Range("N3").Value = -19657.143

# Row 4
$ws.Range("H4").Value = 341.2857
$ws.Range("I4").Value = 176.33333
$ws.Range("J4").Value = 465
$ws.Range("K4").Value = 176.33333
$ws.Range("L4").Value = 465
$ws.Range("M4").Value = -60.33332999999999
$ws.Range("N4").Value = -697

# Row 45
$ws.Range("H45").Value = 1443.6
$ws.Range("I45").Value = 1439.2142
$ws.Range("J45").Value = 1461.1428
$ws.Range("K45").Value = 1439.2142
$ws.Range("L45").Value = 1461.1428
$ws.Range("M45").Value = -1062.2142
$ws.Range("N45").Value = -2215.1428

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 94.23077000000001
$ws.Range("I7").Value = 97.72727
$ws.Range("J7").Value = 75
$ws.Range("K7").Value = 97.72727
$ws.Range("L7").Value = 75
$ws.Range("M7").Value = 15.27273
$ws.Range("N7").Value = -301

# Row 31
$ws.Range("H31").Value = 2297.268
$ws.Range("I31").Value = 1519.6666
$ws.Range("J31").Value = 4536.76
$ws.Range("K31").Value = 1519.6666
$ws.Range("L31").Value = 4536.76
$ws.Range("M31").Value = -1224.6666
$ws.Range("N31").Value = -5126.76

# Row 34
$ws.Range("H34").Value = 2297.268
$ws.Range("I34").Value = 1519.6666
$ws.Range("J34").Value = 4536.76
$ws.Range("K34").Value = 1519.6666
$ws.Range("L34").Value = 4536.76
$ws.Range("M34").Value = -1317.6666
$ws.Range("N34").Value = -4940.76

# Row 58
$ws.Range("H58").Value = 2395257.8
$ws.Range("I58").Value = 4786904.5
$ws.Range("J58").Value = 3611.2104
$ws.Range("K58").Value = 4786904.5
$ws.Range("L58").Value = 3611.2104
$ws.Range("M58").Value = -4786701.5
$ws.Range("N58").Value = -4017.2104

# Row 99
$ws.Range("H99").Value = 3200
$ws.Range("I99").Value = 3240
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 3240
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1742
$ws.Range("N99").Value = -5996

# Row 126
$ws.Range("H126").Value = 3200
$ws.Range("I126").Value = 3240
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 9720
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -7250
$ws.Range("N126").Value = -13940

# Row 136
$ws.Range("H136").Value = 2395257.8
$ws.Range("I136").Value = 4786904.5
$ws.Range("J136").Value = 3611.2104
$ws.Range("K136").Value = 14360713.5
$ws.Range("L136").Value = 10833.6312
$ws.Range("M136").Value = -14358163.5
$ws.Range("N136").Value = -15933.6312

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 95.833336
$ws.Range("J38").Value = 124.583336
$ws.Range("L38").Value = 373.750008
$ws.Range("N38").Value = -1067.750008

# Row 64
$ws.Range("H64").Value = 3475.5
$ws.Range("I64").Value = 1995
$ws.Range("J64").Value = 3969
$ws.Range("K64").Value = 5985
$ws.Range("L64").Value = 11907
$ws.Range("M64").Value = -5715
$ws.Range("N64").Value = -12447

# Row 67
$ws.Range("H67").Value = 3475.5
$ws.Range("I67").Value = 1995
$ws.Range("J67").Value = 3969
$ws.Range("K67").Value = 5985
$ws.Range("L67").Value = 11907
$ws.Range("M67").Value = -5049
$ws.Range("N67").Value = -13779

# Row 70
$ws.Range("H70").Value = 3576.3635
$ws.Range("I70").Value = 1723.5
$ws.Range("J70").Value = 5799.8
$ws.Range("K70").Value = 5170.5
$ws.Range("L70").Value = 17399.4
$ws.Range("M70").Value = -4855.5
$ws.Range("N70").Value = -18029.4

# Row 73
$ws.Range("H73").Value = 3576.3635
$ws.Range("I73").Value = 1723.5
$ws.Range("J73").Value = 5799.8
$ws.Range("K73").Value = 5170.5
$ws.Range("L73").Value = 17399.4
$ws.Range("M73").Value = -4078.5
$ws.Range("N73").Value = -19583.4

# Row 113
$ws.Range("H113").Value = 768.7284
$ws.Range("I113").Value = 781.65
$ws.Range("J113").Value = 731.8095
$ws.Range("K113").Value = 2344.95
$ws.Range("L113").Value = 2195.4285
$ws.Range("M113").Value = -174.9499999999998
$ws.Range("N113").Value = -6535.4285

# Row 127
$ws.Range("H127").Value = 3455.5278
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 3455.5278
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 10366.5834
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -20286.5834

# Row 131
$ws.Range("H131").Value = 522.03
$ws.Range("I131").Value = 285.2586
$ws.Range("J131").Value = 849
$ws.Range("K131").Value = 855.7758
$ws.Range("L131").Value = 2547
$ws.Range("M131").Value = 4184.2242
$ws.Range("N131").Value = -12627

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1866.6666
$ws.Range("I16").Value = 1866.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1866.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1696.6666
$ws.Range("N16").ClearContents()

# Row 22
$ws.Range("H22").Value = 919
$ws.Range("I22").Value = 1450
$ws.Range("J22").Value = 742
$ws.Range("K22").Value = 1450
$ws.Range("L22").Value = 742
$ws.Range("M22").Value = -1155
$ws.Range("N22").Value = -1332

# Row 27
$ws.Range("H27").Value = 919
$ws.Range("I27").Value = 1450
$ws.Range("J27").Value = 742
$ws.Range("K27").Value = 1450
$ws.Range("L27").Value = 742
$ws.Range("M27").Value = -1343
$ws.Range("N27").Value = -956

# Row 46
$ws.Range("H46").Value = 968.75
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 1007.1429
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 1007.1429
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1383.1429

Write-Host "Applied all Pandaemonium_Profits market-data updates"